$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 13. This shifts existing rows 13..124 down to 14..125.
$ws.Rows("13:13").Insert()

# Fill in the new row 13 with the new data record.
$ws.Cells.Item(13, 1).Value = 3
$ws.Cells.Item(13, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(13, 3).Value = "Coquimbo"
$ws.Cells.Item(13, 4).Value = 44490
$ws.Cells.Item(13, 5).Value = 5
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100101
$ws.Cells.Item(13, 8).Value = "Berries"
$ws.Cells.Item(13, 9).Value = 100101001
$ws.Cells.Item(13, 10).Value = "Arándano (blue)"
$ws.Cells.Item(13, 11).Value = "Sin especificar"
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 116
$ws.Cells.Item(13, 14).Value = 10000
$ws.Cells.Item(13, 15).Value = 11000
$ws.Cells.Item(13, 16).Value = 10483
$ws.Cells.Item(13, 17).Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(13, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(13, 19).Value = 6989
$ws.Cells.Item(13, 20).Value = 1.5
